# Incremento Clase Corredor en la planilla de métricas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 of the "Desarrollo y correctivos" increment table.
$ws.Range("C18").Value = "Clase Corredor"
$ws.Range("F18").Value = 50
$ws.Range("G18").Value = 10.0 / 1440.0
$ws.Range("H18").Value = 0.86805555555555547
$ws.Range("I18").Value = 0.87361111111111101
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 45

# Move the active selection, matching where the user clicked next.
$ws.Range("F19").Select()
